$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 97
$ws.Range("C97").Value = 1469681.5544
$ws.Range("D97").Value = 1447357.2265000001
$ws.Range("E97").Value = 1960.5377000000001
$ws.Range("F97").Value = 589939.82200000004
$ws.Range("G97").Value = 176730.3
$ws.Range("H97").Value = 152891.6
$ws.Range("I97").Value = 336207.31
$ws.Range("J97").Value = 17588.062399999999
$ws.Range("K97").Value = 450626.97869999998
$ws.Range("L97").Value = 683239.71470000001
$ws.Range("M97").Value = 419208.54220000003
$ws.Range("N97").Value = 73123.818700000003
$ws.Range("O97").Value = 62472.796600000001
$ws.Range("P97").Value = 8916.6134999999995
$ws.Range("Q97").Value = 10864.0401
$ws.Range("R97").Value = 141345.74
$ws.Range("S97").Value = 134349.39600000001
$ws.Range("T97").Value = 117970.283
$ws.Range("U97").Value = 21450.388999999999
$ws.Range("V97").Value = 5076.28
$ws.Range("W97").Value = 687.23200285155394
$ws.Range("X97").Value = 498330.75400000002
$ws.Range("Y97").Value = 23482.256000000001
$ws.Range("AC97").Value = 1955.63258
$ws.Range("AD97").Value = 826.12345000000005
$ws.Range("AE97").Value = 464.53843000000001

# Row 98
$ws.Range("C98").Value = 630368
$ws.Range("D98").Value = 602983
$ws.Range("E98").Value = 1715.6668999999999
$ws.Range("F98").Value = 191387.4
$ws.Range("G98").Value = 92343.4
$ws.Range("H98").Value = 82169.899999999994
$ws.Range("I98").Value = 178975.14
$ws.Range("J98").Value = 15057.7844
$ws.Range("K98").Value = 257334.97940000001
$ws.Range("L98").Value = 652647.12459999998
$ws.Range("M98").Value = 229450.60680000001
$ws.Range("N98").Value = 38978.090900000003
$ws.Range("O98").Value = 11846.3465
$ws.Range("P98").Value = 33918.995699999999
$ws.Range("Q98").Value = 29131.8842
$ws.Range("R98").Value = 183786.71900000001
$ws.Range("S98").Value = 179113.20300000001
$ws.Range("T98").Value = 139259.41500000001
$ws.Range("U98").Value = 41898.345000000001
$ws.Range("V98").Value = 2632.5
$ws.Range("W98").Value = 491.21483270033099
$ws.Range("X98").Value = 196069.80900000001
$ws.Range("Y98").Value = 6248.3069999999998
$ws.Range("Z98").Value = 6499.3127599999998
$ws.Range("AA98").Value = 744.69
$ws.Range("AC98").Value = 7213.0374300000003
$ws.Range("AD98").Value = 2525.4867399999998
$ws.Range("AE98").Value = 976.24905999999999
